$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-7 (A=Sending cluster, B=Ligand, C=Receptor, D=Target cluster,
# E..T numeric columns). Rows 2-4 get updated in place; rows 5-7 are newly added.
$rows = @(
    ,@("ECs", "FAPs", @(3, 1, 91.20678233333332, 273.620347, 0.6532435006323181, 0.6532435006323182, 3, 1, 2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 182.6337682414801, 1643.703914173321, 0.5341090210140617, 0.5341090210140617))
    ,@("ECs", "sCs", @(3, 1, 91.20678233333332, 273.620347, 0.6532435006323181, 0.6532435006323182, 2, 0.6666666666666666, 0.446644, 1.339932, 0.1823737695100498, 0.1823737695100498, 40.73696208848934, 366.632658796404, 0.1191344796182564, 0.1191344796182565))
    ,@("FAPs", "FAPs", @(3, 1, 27.97197233333334, 83.91591700000001, 0.2003415607825798, 0.2003415607825798, 3, 1, 2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 56.01147833187012, 504.1033049868311, 0.1638045151531339, 0.1638045151531339))
    ,@("FAPs", "sCs", @(3, 1, 27.97197233333334, 83.91591700000001, 0.2003415607825798, 0.2003415607825798, 2, 0.6666666666666666, 0.446644, 1.339932, 0.1823737695100498, 0.1823737695100498, 12.49351361084934, 112.441622497644, 0.03653704562944583, 0.03653704562944584))
    ,@("sCs", "FAPs", @(3, 1, 20.442661, 61.327983, 0.146414938585102, 0.146414938585102, 3, 1, 2.002414333333333, 6.007243, 0.8176262304899502, 0.8176262304899502, 40.93467739787433, 368.412096580869, 0.1197126943227545, 0.1197126943227545))
    ,@("sCs", "sCs", @(3, 1, 20.442661, 61.327983, 0.146414938585102, 0.146414938585102, 2, 0.6666666666666666, 0.446644, 1.339932, 0.1823737695100498, 0.1823737695100498, 9.130591879684001, 82.17532691715601, 0.02670224426234749, 0.02670224426234749))
)

$r = 2
foreach ($row in $rows) {
    $a = $row[0]
    $d = $row[1]
    $vals = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = "Timp3"
    $ws.Cells.Item($r, 3).Value = "Agtr2"
    $ws.Cells.Item($r, 4).Value = $d
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $vals[$i]
    }
    $r++
}
